$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46030
$ws.Range("B2").Value = 84.79000000000001
$ws.Range("C2").Value = 83.58
$ws.Range("D2").Value = 82.18000000000001
$ws.Range("E2").Value = 79.06
$ws.Range("F2").Value = 74.47
$ws.Range("G2").Value = 77.98
$ws.Range("H2").Value = 81.05
$ws.Range("I2").Value = 95.75
$ws.Range("J2").Value = 107.36
$ws.Range("K2").Value = 95.66
$ws.Range("L2").Value = 82.56
$ws.Range("M2").Value = 68.84
$ws.Range("N2").Value = 65.73999999999999
$ws.Range("O2").Value = 63.55
$ws.Range("P2").Value = 60.08
$ws.Range("Q2").Value = 63.4
$ws.Range("R2").Value = 76.86
$ws.Range("S2").Value = 85.84999999999999
$ws.Range("T2").Value = 97.87
$ws.Range("U2").Value = 98.73999999999999
$ws.Range("V2").Value = 93.63
$ws.Range("W2").Value = 86.77
$ws.Range("X2").Value = 82.98
$ws.Range("Y2").Value = 74.34999999999999
$ws.Range("Z2").Value = 81.8
$ws.Range("AB2").Value = 89.83
$ws.Range("AC2").Value = "8h-10h"
$ws.Range("AD2").Value = 101.51
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 98.3
$ws.Range("AG2").Value = "3h-23h"
